$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with new recipe names, including new column G
$ws.Range("B1").Value = "Recipe Lomo Saltado with beef"
$ws.Range("C1").Value = "Recipe Lomo Saltado with chicken"
$ws.Range("D1").Value = "Recipe Lomo Saltado with chickpeas"
$ws.Range("E1").Value = "Recipe Lomo Saltado with mushrooms"
$ws.Range("F1").Value = "Recipe Lomo Saltado with salmon"
$ws.Range("G1").Value = "Recipe Lomo Saltado with tofu"

# Copy the header style from F1 to the newly introduced G1 header cell
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122) # xlPasteFormats

# Row 2 - Acidification
$ws.Range("B2").Value = 0.107303026
$ws.Range("C2").Value = 0.050473526
$ws.Range("D2").Value = 0.029880226
$ws.Range("E2").Value = 0.023592126
$ws.Range("F2").Value = 0.032551076
$ws.Range("G2").Value = 0.024019826

# Row 3 - Eutrophication
$ws.Range("B3").Value = 0.103184908
$ws.Range("C3").Value = 0.035890408
$ws.Range("D3").Value = 0.029265608
$ws.Range("E3").Value = 0.023282358
$ws.Range("F3").Value = 0.064327908
$ws.Range("G3").Value = 0.024952208

# Row 4 - Freshwater Withdrawals (FW)
$ws.Range("B4").Value = 1432.1124
$ws.Range("C4").Value = 1218.7174
$ws.Range("D4").Value = 1201.4274
$ws.Range("E4").Value = 1073.5724
$ws.Range("F4").Value = 1634.1324
$ws.Range("G4").Value = 1099.0524

# Row 5 - Scarcity-Weighted FW
$ws.Range("B5").Value = 40831.37280000001
$ws.Range("C5").Value = 35218.94780000002
$ws.Range("D5").Value = 39697.05780000002
$ws.Range("E5").Value = 33003.0978
$ws.Range("F5").Value = 37431.61280000001
$ws.Range("G5").Value = 33379.38280000001

# Row 6 - GHG emissions
$ws.Range("B6").Value = 31.6705648
$ws.Range("C6").Value = 7.3380748
$ws.Range("D6").Value = 5.2987648
$ws.Range("E6").Value = 4.773239800000001
$ws.Range("F6").Value = 6.924479800000001
$ws.Range("G6").Value = 5.872519799999999

# Row 7 - Land use
$ws.Range("B7").Value = 93.2142344
$ws.Range("C7").Value = 8.504154400000001
$ws.Range("D7").Value = 11.2933044
$ws.Range("E7").Value = 5.2372544
$ws.Range("F7").Value = 6.633194400000001
$ws.Range("G7").Value = 6.531274400000001
